$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format to Text for price cells so numeric-looking strings
# (e.g. "7.30", "44.388.37") are preserved exactly, matching the source data
# which stores these as literal text, not numbers.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '44.388.37'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '2.250.71'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '307.72'
$ws.Range('E5').Value = '  -1.97%  '
$ws.Range('D6').Value = '96.80'
$ws.Range('E6').Value = '  -1.50%  '
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '7.30'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '2.592.85'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').Value = '2.245.55'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '0.837'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '13.67'
$ws.Range('E17').Value = '  -2.28%  '
$ws.Range('D18').Value = '44.179.38'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('E20').Value = '  -7.03%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').Value = '65.78'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').Value = '238.29'
$ws.Range('E23').Value = '  +1.15%  '
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '38.94'
$ws.Range('E27').Value = '  +6.69%  '
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('D29').Value = '2.21'
$ws.Range('E29').Value = '  +3.32%  '
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').Value = '152.50'
$ws.Range('E32').Value = '  -2.69%  '
$ws.Range('D33').Value = '0.0801'
$ws.Range('E33').Value = '  -3.67%  '
$ws.Range('E34').Value = '  -1.83%  '
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('E36').Value = '  +2.72%  '
$ws.Range('E37').Value = '  -1.98%  '
$ws.Range('D38').Value = '1.78'
$ws.Range('E38').Value = '  -6.45%  '
$ws.Range('D39').Value = '3.62'
$ws.Range('E39').Value = '  +2.43%  '
$ws.Range('D40').Value = '14.67'
$ws.Range('E40').Value = '  -6.63%  '
$ws.Range('D41').Value = '3.87'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('E42').Value = '  -1.77%  '
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').Value = '1.754.14'
$ws.Range('E44').Value = '  +2.97%  '
$ws.Range('D45').Value = '83.30'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('D47').Value = '15.32'
$ws.Range('E47').Value = '  +12.68%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '4.98'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '100.56'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').Value = '8.18'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('E51').Value = '  -1.95%  '
